$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "current" sheet: add a row for 2024-05-27 (doctor Florentiy_Pavlov, id 6)
# ---------------------------------------------------------------------------
$current = $wb.Worksheets.Item("current")

$current.Cells.Item(7, 1).NumberFormat = "@"
$current.Cells.Item(7, 1).Value = "2024-05-27"
$current.Cells.Item(7, 1).ClearFormats()

$current.Cells.Item(7, 2).Value = "Florentiy_Pavlov"
$current.Cells.Item(7, 3).Value = 0

$current.Cells.Item(7, 4).NumberFormat = "@"
$current.Cells.Item(7, 4).Value = "6"
$current.Cells.Item(7, 4).ClearFormats()

# ---------------------------------------------------------------------------
# 2) Remember which sheet is currently active/selected, so we can restore it
#    after inserting the new sheet at the end of the tab strip.
# ---------------------------------------------------------------------------
$previouslyActive = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 3) Add a brand-new sheet "2024-05-27" after the last existing sheet
#    ("2024-05-26") and fill it with today's patient log.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2024-05-27"

$newSheet.Cells.Item(1, 1).Value = "Время"
$newSheet.Cells.Item(1, 2).Value = "ФИО пациента"
$newSheet.Cells.Item(1, 3).Value = "М\Ж\Р"

$rows = @(
    @(45439.75094185185, "jbnkjno", "Ж"),
    @(45439.83571817129, "sjfn'pdfi", "Ж"),
    @(45439.83577740741, "fkmj'fdpk", "Ж"),
    @(45439.83584238426, "fkmj'fdpk", "Р"),
    @(45439.83796740741, "gdxjchvbk", "M"),
    @(45439.83801888889, "fxckgvbhjk", "M"),
    @(45439.83807405092, "sjfn'pdfi", "Р"),
    @(45439.83813391204, "gvkjhkj", "Ж"),
    @(45439.83820114583, "fkmj'fdpk", "Р"),
    @(45439.83825736111, "jhblk", "M"),
    @(45439.8383075, "aedf", "Ж"),
    @(45439.83836677083, "aedf", "M"),
    @(45439.84088684028, "fkmj'fdpk", "Ж"),
    @(45439.85151035563, "sjfn'pdfi", "M")
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Restore the tab that was active before we inserted the new sheet.
$previouslyActive.Activate()

# ---------------------------------------------------------------------------
# 4) "2024-05-26" sheet: A5 timestamp got a tiny precision bump on save.
# ---------------------------------------------------------------------------
$may26 = $wb.Worksheets.Item("2024-05-26")
$may26.Cells.Item(5, 1).Value = 45438.84244952547

Write-Output "ok"
